# Refresh cryptos list: update Price (D) and Volume(1h) (E) columns with latest values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.744.27'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.47%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.910.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.14%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.82%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.67%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.555'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.20%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.916.85'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.29%  '

$ws.Range("E10").Value = '  -4.84%  '

$ws.Range("E11").Value = '  -1.95%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.361'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.59%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.417.82'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.11%  '

$ws.Range("E14").Value = '  +2.35%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.705.87'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.66%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.83'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.914.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.30%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000142'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.95%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.73'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '361.64'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.53%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.54%  '

$ws.Range("E23").Value = '  -0.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.67'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.19%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.456'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.22%  '

$ws.Range("E27").Value = '  -3.86%  '

$ws.Range("E28").Value = '  +0.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.89'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0851'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.93%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.02%  '

$ws.Range("E32").Value = '  -1.68%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.86'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.37'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.79%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.60'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.19%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.01'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.21'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.94'
$ws.Range("D39").Style = "Normal"

$ws.Range("E40").Value = '  -4.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.73'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.299.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.78%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.651'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.50%  '

$ws.Range("E44").Value = '  -1.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.84%  '

$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.36%  '

$ws.Range("E48").Value = '  -3.76%  '

$ws.Range("E49").Value = '  -1.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0924'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '250.32'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.24%  '
